# Daily attendance processing - 2025-12-19 17:27:32
# Reorders the comma-separated "Recorded By" values in column G so that any
# "System" / "system" entries are moved to the end of the list (in reverse
# order among themselves), while the remaining (non-system) entries keep
# their original relative order at the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column G (header is in row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val -split ",\s*"

    $nonSystem = @()
    $systemParts = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemParts += $trimmed
        } else {
            $nonSystem += $trimmed
        }
    }

    # reverse the order of the system/System entries
    $revSystem = @()
    for ($i = $systemParts.Count - 1; $i -ge 0; $i--) {
        $revSystem += $systemParts[$i]
    }

    $newParts = $nonSystem + $revSystem
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
